$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 350, shifting rows 350:365 down to 351:366.
$ws.Rows(350).Insert()

# Populate the newly inserted row 350 with the new record's data.
$ws.Range("A350").Value = 7
$ws.Range("B350").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C350").Value = "Ñuble"
$ws.Range("D350").Value2 = 44939
$ws.Range("E350").Value = 16
$ws.Range("F350").Value = 100114013
$ws.Range("G350").Value = "Zanahoria"
$ws.Range("H350").Value = "Sin especificar"
$ws.Range("I350").Value = "Primera"
$ws.Range("J350").Value = 160
$ws.Range("K350").Value = 8500
$ws.Range("L350").Value = 9000
$ws.Range("M350").Value = 8750
$ws.Range("N350").Value = "$/saco 20 kilos"
$ws.Range("O350").Value = "Región de Ñuble"
$ws.Range("P350").Value = 438
$ws.Range("Q350").Value = 20
$ws.Range("R350").Value = "Hortaliza"
